$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "243.31"
Set-TextValue "D4" "5.407"
Set-TextValue "D5" "0.05988"
Set-TextValue "D6" "3.429"
Set-TextValue "D7" "6.491"
Set-TextValue "D8" "0.8085"
Set-TextValue "D9" "0.9240"
Set-TextValue "D10" "0.1428"
Set-TextValue "D11" "0.07413"
Set-TextValue "D12" "0.03261"
Set-TextValue "D13" "0.03071"
Set-TextValue "D15" "3.859"
Set-TextValue "D16" "0.001574"
Set-TextValue "D17" "0.04694"
Set-TextValue "D18" "0.0005899"
$ws.Range("E18").Value = "17OneONEWorstin24h"
Set-TextValue "D19" "0.005857"
Set-TextValue "D20" "0.001265"
Set-TextValue "D23" "3.571"
Set-TextValue "D24" "2.179"
Set-TextValue "D26" "0.1331"
Set-TextValue "D27" "0.0002339"
Set-TextValue "D40" "0.03970"
Set-TextValue "D41" "0.006390"
Set-TextValue "D43" "0.1078"
Set-TextValue "D44" "0.008602"
Set-TextValue "D45" "0.00005094"
Set-TextValue "D47" "0.6499"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
Set-TextValue "D48" "0.002452"
